$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "system, backup@backdoor.com, System"
$ws.Range("G3").Value = "System, dnasr281@gmail.com"
$ws.Range("G4").Value = "backup@backdoor.com, System"
$ws.Range("G5").Value = "backup@backdoor.com, System"
$ws.Range("G6").Value = "System, dnasr281@gmail.com"
$ws.Range("G11").Value = "System, dnasr281@gmail.com"
$ws.Range("G12").Value = "System, dnasr281@gmail.com"
$ws.Range("G13").Value = "System, dnasr281@gmail.com"
$ws.Range("G14").Value = "System, dnasr281@gmail.com"
$ws.Range("G15").Value = "System, dnasr281@gmail.com"
$ws.Range("G29").Value = "system, backup@backdoor.com, System"
$ws.Range("G30").Value = "System, dnasr281@gmail.com"
$ws.Range("G32").Value = "backup@backdoor.com, System"
$ws.Range("G33").Value = "System, dnasr281@gmail.com"
$ws.Range("G38").Value = "System, dnasr281@gmail.com"
$ws.Range("G39").Value = "System, dnasr281@gmail.com"
$ws.Range("G40").Value = "System, dnasr281@gmail.com"
$ws.Range("G41").Value = "System, dnasr281@gmail.com"
$ws.Range("G42").Value = "System, dnasr281@gmail.com"
$ws.Range("G56").Value = "system, backup@backdoor.com, System"
$ws.Range("G57").Value = "System, dnasr281@gmail.com"
$ws.Range("G58").Value = "backup@backdoor.com, System"
$ws.Range("G59").Value = "backup@backdoor.com, System"
$ws.Range("G60").Value = "System, dnasr281@gmail.com"
$ws.Range("G65").Value = "System, dnasr281@gmail.com"
$ws.Range("G66").Value = "System, dnasr281@gmail.com"
$ws.Range("G67").Value = "System, dnasr281@gmail.com"
$ws.Range("G68").Value = "System, dnasr281@gmail.com"
$ws.Range("G69").Value = "System, dnasr281@gmail.com"
$ws.Range("G84").Value = "backup@backdoor.com, System"
$ws.Range("G85").Value = "backup@backdoor.com, System"
$ws.Range("G86").Value = "System, dnasr281@gmail.com"
$ws.Range("G89").Value = "System, dnasr281@gmail.com"
$ws.Range("G90").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G93").Value = "System, dnasr281@gmail.com"
$ws.Range("G110").Value = "backup@backdoor.com, System"
$ws.Range("G111").Value = "backup@backdoor.com, System"
$ws.Range("G112").Value = "System, dnasr281@gmail.com"
$ws.Range("G115").Value = "System, dnasr281@gmail.com"
$ws.Range("G116").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G119").Value = "System, dnasr281@gmail.com"
$ws.Range("G136").Value = "backup@backdoor.com, System"
$ws.Range("G137").Value = "backup@backdoor.com, System"
$ws.Range("G138").Value = "System, dnasr281@gmail.com"
$ws.Range("G141").Value = "System, dnasr281@gmail.com"
$ws.Range("G142").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G145").Value = "System, dnasr281@gmail.com"
